$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "202.173.124.249"
$ws.Range("B4").Value = 28.5212672
$ws.Range("C4").Value = 77.2243456
$ws.Range("D4").Value = 735485.4806669627
$ws.Range("E4").Value = "Mozilla/5.0 (Windows NT 10.0; Win64; x64) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Safari/537.36"
$ws.Range("F4").Value = "Win32"
$ws.Range("G4").Value = "2025-06-21T04:19:34.989Z"

$ws.Range("A5").Value = "202.173.124.249"
$ws.Range("B5").Value = 28.3621566
$ws.Range("C5").Value = 77.2827572
$ws.Range("D5").Value = 15.079999923706055
$ws.Range("E5").Value = "Mozilla/5.0 (Linux; Android 10; K) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Mobile Safari/537.36"
$ws.Range("F5").Value = "Linux armv81"
$ws.Range("G5").Value = "2025-06-21T04:20:03.844Z"
